$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Price" column (D) values - force text type to match the
# original inline-string cells, then reset the style so no stray
# number-format style is left behind.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '71.296.02'
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.814.70'
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '702.85'
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '171.16'
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.813.93'
$ws.Range("D7").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.489'
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.47'
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000251'
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.72'
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.454.37'
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.824.00'
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '71.380.56'
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.24'
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.54'
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '510.94'
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.47'
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '83.79'
$ws.Range("D24").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.78'
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.960.94'
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.34'
$ws.Range("D28").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.37'
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '29.20'
$ws.Range("D34").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '9.31'
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.779.28'
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.999'
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.68'
$ws.Range("D39").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.23'
$ws.Range("D44").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '166.49'
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '50.07'
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '433.89'
$ws.Range("D48").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '30.69'
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.39'
$ws.Range("D51").Style = "Normal"

# Update "Volume(1h)" column (E) values.
$ws.Range("E2").Value = '  +0.17%  '
$ws.Range("E3").Value = '  -0.64%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("E5").Value = '  -1.58%  '
$ws.Range("E6").Value = '  -1.13%  '
$ws.Range("E7").Value = '  -0.62%  '
$ws.Range("E8").Value = '  +0.09%  '
$ws.Range("E9").Value = '  -0.28%  '
$ws.Range("E10").Value = '  -2.01%  '
$ws.Range("E11").Value = '  +6.28%  '
$ws.Range("E12").Value = '  +1.50%  '
$ws.Range("E13").Value = '  -2.35%  '
$ws.Range("E14").Value = '  -0.31%  '
$ws.Range("E15").Value = '  -0.71%  '
$ws.Range("E16").Value = '  +1.39%  '
$ws.Range("E17").Value = '  +0.27%  '
$ws.Range("E18").Value = '  +0.13%  '
$ws.Range("E19").Value = '  +0.76%  '
$ws.Range("E20").Value = '  +0.04%  '
$ws.Range("E21").Value = '  +3.25%  '
$ws.Range("E22").Value = '  -2.47%  '
$ws.Range("E23").Value = '  -1.58%  '
$ws.Range("E24").Value = '  -1.67%  '
$ws.Range("E25").Value = '  -3.16%  '
$ws.Range("E26").Value = '  +4.87%  '
$ws.Range("E27").Value = '  -0.72%  '
$ws.Range("E28").Value = '  -3.13%  '
$ws.Range("E29").Value = '  +0.03%  '
$ws.Range("E30").Value = '  -4.58%  '
$ws.Range("E31").Value = '  -5.68%  '
$ws.Range("E32").Value = '  +0.88%  '
$ws.Range("E33").Value = '  -1.76%  '
$ws.Range("E34").Value = '  -0.77%  '
$ws.Range("E35").Value = '  -5.33%  '
$ws.Range("E36").Value = '  +0.79%  '
$ws.Range("E37").Value = '  -0.63%  '
$ws.Range("E38").Value = '  +0.20%  '
$ws.Range("E39").Value = '  +11.14%  '
$ws.Range("E40").Value = '  -2.11%  '
$ws.Range("E41").Value = '  +5.29%  '
$ws.Range("E42").Value = '  -2.31%  '
$ws.Range("E43").Value = '  +0.01%  '
$ws.Range("E44").Value = '  -4.02%  '
$ws.Range("E45").Value = '  -0.06%  '
$ws.Range("E46").Value = '  +1.69%  '
$ws.Range("E47").Value = '  +2.44%  '
$ws.Range("E48").Value = '  +2.24%  '
$ws.Range("E49").Value = '  -5.67%  '
$ws.Range("E50").Value = '  +9.04%  '
$ws.Range("E51").Value = '  +0.02%  '
